$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.2845560014247894
$ws.Range("B1").Value = 0.8495793342590332
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 1.115927577018738
